$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Activity #1): keep B2's string slot but update text; change C2 text
$ws.Range("B2").Value = "Run auto register script with real data and fix some bugs"
$ws.Range("C2").Value = "Didn't find the root cause why the browser opened more than 8 while I set Pool size to 8"
$ws.Rows(2).RowHeight = 27

# Row 3 (Result #1): fill in B3/C3 with what used to live in row 6
$ws.Range("B3").Value = "learn Python data types: string"
$ws.Range("C3").Value = "learned by Think like a computer scientist"

# Row 6 (Activity #2): clear out text, restore default row height
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = ""
$ws.Rows(6).AutoFit()

# Update the visible selection to A1:C3
$ws.Range("A1:C3").Select()
